$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 32
$ws.Cells.Item(32, 8).Value = 8299.25
$ws.Cells.Item(32, 10).Value = 3248
$ws.Cells.Item(32, 12).Value = 3248
$ws.Cells.Item(32, 14).Value = -3900

# ALC row 43
$ws.Cells.Item(43, 8).Value = 8500
$ws.Cells.Item(43, 9).Value = 8500
$ws.Cells.Item(43, 10).Value = 8500
$ws.Cells.Item(43, 11).Value = 8500
$ws.Cells.Item(43, 12).Value = 8500
$ws.Cells.Item(43, 13).Value = -8431
$ws.Cells.Item(43, 14).Value = -8638

# ALC row 110
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).ClearContents()
$ws.Cells.Item(110, 14).Value = 0

# ALC row 129
$ws.Cells.Item(129, 8).Value = 2084246.2
$ws.Cells.Item(129, 9).Value = 744.7646999999999
$ws.Cells.Item(129, 11).Value = 2234.2941
$ws.Cells.Item(129, 13).Value = 2765.7059

# ALC row 137
$ws.Cells.Item(137, 8).Value = 3813.4827
$ws.Cells.Item(137, 9).Value = 2761.5
$ws.Cells.Item(137, 10).Value = 4214.2383
$ws.Cells.Item(137, 11).Value = 8284.5
$ws.Cells.Item(137, 12).Value = 12642.7149
$ws.Cells.Item(137, 13).Value = -5734.5
$ws.Cells.Item(137, 14).Value = -17742.7149

# ALC row 138
$ws.Cells.Item(138, 8).Value = 3799.4167
$ws.Cells.Item(138, 9).Value = 1391.1578
$ws.Cells.Item(138, 10).Value = 5377.241
$ws.Cells.Item(138, 11).Value = 4173.4734
$ws.Cells.Item(138, 12).Value = 16131.723
$ws.Cells.Item(138, 13).Value = 966.5266000000001
$ws.Cells.Item(138, 14).Value = -26411.723

# ALC row 141
$ws.Cells.Item(141, 8).Value = 2620.1462
$ws.Cells.Item(141, 9).Value = 2473.525
$ws.Cells.Item(141, 10).Value = 8485
$ws.Cells.Item(141, 11).Value = 7420.575000000001
$ws.Cells.Item(141, 12).Value = 25455
$ws.Cells.Item(141, 13).Value = -2240.575000000001
$ws.Cells.Item(141, 14).Value = -35815

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Cells.Item(32, 8).Value = 5275.3335
$ws.Cells.Item(32, 9).Value = 4997.1
$ws.Cells.Item(32, 10).Value = 15013.5
$ws.Cells.Item(32, 11).Value = 4997.1
$ws.Cells.Item(32, 12).Value = 15013.5
$ws.Cells.Item(32, 13).Value = -4710.1
$ws.Cells.Item(32, 14).Value = -15587.5

# ARM row 61
$ws.Cells.Item(61, 8).Value = 8932675
$ws.Cells.Item(61, 9).Value = 9095051
$ws.Cells.Item(61, 10).Value = 1993
$ws.Cells.Item(61, 11).Value = 9095051
$ws.Cells.Item(61, 12).Value = 1993
$ws.Cells.Item(61, 13).Value = -9094839
$ws.Cells.Item(61, 14).Value = -2417

# ARM row 74
$ws.Cells.Item(74, 8).Value = 3336.9348
$ws.Cells.Item(74, 9).Value = 2198.6487
$ws.Cells.Item(74, 10).Value = 8016.5557
$ws.Cells.Item(74, 11).Value = 2198.6487
$ws.Cells.Item(74, 12).Value = 8016.5557
$ws.Cells.Item(74, 13).Value = -1324.6487
$ws.Cells.Item(74, 14).Value = -9764.555700000001

# ARM row 77
$ws.Cells.Item(77, 8).Value = 3336.9348
$ws.Cells.Item(77, 9).Value = 2198.6487
$ws.Cells.Item(77, 10).Value = 8016.5557
$ws.Cells.Item(77, 11).Value = 10993.2435
$ws.Cells.Item(77, 12).Value = 40082.7785
$ws.Cells.Item(77, 13).Value = -6625.2435
$ws.Cells.Item(77, 14).Value = -48818.7785

# ARM row 97
$ws.Cells.Item(97, 8).Value = 1504.3214
$ws.Cells.Item(97, 9).Value = 1179.2174
$ws.Cells.Item(97, 10).Value = 2999.8
$ws.Cells.Item(97, 11).Value = 1179.2174
$ws.Cells.Item(97, 12).Value = 2999.8
$ws.Cells.Item(97, 13).Value = -683.2174
$ws.Cells.Item(97, 14).Value = -3991.8

# ARM row 106
$ws.Cells.Item(106, 8).Value = 42685
$ws.Cells.Item(106, 10).Value = 42685
$ws.Cells.Item(106, 12).Value = 42685
$ws.Cells.Item(106, 14).Value = -45209

# ARM row 124
$ws.Cells.Item(124, 8).Value = 30429
$ws.Cells.Item(124, 10).Value = 30429
$ws.Cells.Item(124, 12).Value = 30429
$ws.Cells.Item(124, 14).Value = -40249

# ARM row 132
$ws.Cells.Item(132, 8).Value = 2596.1128
$ws.Cells.Item(132, 9).Value = 2083.9807
$ws.Cells.Item(132, 11).Value = 6251.9421
$ws.Cells.Item(132, 13).Value = -3721.9421

# ARM row 136
$ws.Cells.Item(136, 8).Value = 8932675
$ws.Cells.Item(136, 9).Value = 9095051
$ws.Cells.Item(136, 10).Value = 1993
$ws.Cells.Item(136, 11).Value = 27285153
$ws.Cells.Item(136, 12).Value = 5979
$ws.Cells.Item(136, 13).Value = -27282603
$ws.Cells.Item(136, 14).Value = -11079

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20
$ws.Cells.Item(20, 8).Value = 2404
$ws.Cells.Item(20, 9).Value = 1800.6428
$ws.Cells.Item(20, 10).Value = 3397.7646
$ws.Cells.Item(20, 11).Value = 1800.6428
$ws.Cells.Item(20, 12).Value = 3397.7646
$ws.Cells.Item(20, 13).Value = -1553.6428
$ws.Cells.Item(20, 14).Value = -3891.7646

# BSM row 86
$ws.Cells.Item(86, 8).Value = 1202.6177
$ws.Cells.Item(86, 9).Value = 1120.0385
$ws.Cells.Item(86, 10).Value = 1471
$ws.Cells.Item(86, 11).Value = 1120.0385
$ws.Cells.Item(86, 12).Value = 1471
$ws.Cells.Item(86, 13).Value = 2.961499999999887
$ws.Cells.Item(86, 14).Value = -3717

# BSM row 89
$ws.Cells.Item(89, 8).Value = 1202.6177
$ws.Cells.Item(89, 9).Value = 1120.0385
$ws.Cells.Item(89, 10).Value = 1471
$ws.Cells.Item(89, 11).Value = 5600.192500000001
$ws.Cells.Item(89, 12).Value = 7355
$ws.Cells.Item(89, 13).Value = 15.80749999999898
$ws.Cells.Item(89, 14).Value = -18587

# BSM row 128
$ws.Cells.Item(128, 8).Value = 3000
$ws.Cells.Item(128, 9).Value = 3000
$ws.Cells.Item(128, 11).Value = 9000
$ws.Cells.Item(128, 13).Value = -6510

# BSM row 132
$ws.Cells.Item(132, 8).Value = 88070
$ws.Cells.Item(132, 10).Value = 88070
$ws.Cells.Item(132, 12).Value = 88070
$ws.Cells.Item(132, 14).Value = -98190

# BSM row 134
$ws.Cells.Item(134, 8).Value = 3345.2468
$ws.Cells.Item(134, 9).Value = 3217.0154
$ws.Cells.Item(134, 10).Value = 4039.8333
$ws.Cells.Item(134, 11).Value = 9651.046200000001
$ws.Cells.Item(134, 12).Value = 12119.4999
$ws.Cells.Item(134, 13).Value = -7116.046200000001
$ws.Cells.Item(134, 14).Value = -17189.4999

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16
$ws.Cells.Item(16, 8).Value = 3749.9
$ws.Cells.Item(16, 9).Value = 3187.5
$ws.Cells.Item(16, 10).Value = 5999.5
$ws.Cells.Item(16, 11).Value = 3187.5
$ws.Cells.Item(16, 12).Value = 5999.5
$ws.Cells.Item(16, 13).Value = -2900.5
$ws.Cells.Item(16, 14).Value = -6573.5

# CRP row 23
$ws.Cells.Item(23, 8).Value = 66330.75
$ws.Cells.Item(23, 9).Value = 66329.5
$ws.Cells.Item(23, 10).Value = 66332
$ws.Cells.Item(23, 11).Value = 66329.5
$ws.Cells.Item(23, 12).Value = 66332
$ws.Cells.Item(23, 13).Value = -66089.5
$ws.Cells.Item(23, 14).Value = -66812

# CRP row 27
$ws.Cells.Item(27, 8).Value = 66330.75
$ws.Cells.Item(27, 9).Value = 66329.5
$ws.Cells.Item(27, 10).Value = 66332
$ws.Cells.Item(27, 11).Value = 66329.5
$ws.Cells.Item(27, 12).Value = 66332
$ws.Cells.Item(27, 13).Value = -66137.5
$ws.Cells.Item(27, 14).Value = -66716

# CRP row 36
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 13).ClearContents()

# CRP row 40
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 13).ClearContents()

# CRP row 64
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).ClearContents()
$ws.Cells.Item(64, 14).Value = 0

# CRP row 67
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).ClearContents()
$ws.Cells.Item(67, 14).Value = 0

# CRP row 113
$ws.Cells.Item(113, 8).Value = 3749.9
$ws.Cells.Item(113, 9).Value = 3187.5
$ws.Cells.Item(113, 10).Value = 5999.5
$ws.Cells.Item(113, 11).Value = 3187.5
$ws.Cells.Item(113, 12).Value = 5999.5
$ws.Cells.Item(113, 13).Value = -1017.5
$ws.Cells.Item(113, 14).Value = -10339.5

# CRP row 132
$ws.Cells.Item(132, 8).Value = 3537.5789
$ws.Cells.Item(132, 9).Value = 2264.0625
$ws.Cells.Item(132, 10).Value = 10329.667
$ws.Cells.Item(132, 11).Value = 6792.1875
$ws.Cells.Item(132, 12).Value = 30989.001
$ws.Cells.Item(132, 13).Value = -4262.1875
$ws.Cells.Item(132, 14).Value = -36049.001

$ws = $wb.Worksheets.Item("GSM")
# GSM row 132
$ws.Cells.Item(132, 8).Value = 5534.8
$ws.Cells.Item(132, 9).Value = 5510.4165
$ws.Cells.Item(132, 11).Value = 16531.2495
$ws.Cells.Item(132, 13).Value = -14001.2495

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Cells.Item(22, 8).Value = 1431.0769
$ws.Cells.Item(22, 9).Value = 1732.7142
$ws.Cells.Item(22, 10).Value = 1079.1666
$ws.Cells.Item(22, 11).Value = 1732.7142
$ws.Cells.Item(22, 12).Value = 1079.1666
$ws.Cells.Item(22, 13).Value = -1437.7142
$ws.Cells.Item(22, 14).Value = -1669.1666

# LTW row 27
$ws.Cells.Item(27, 8).Value = 1431.0769
$ws.Cells.Item(27, 9).Value = 1732.7142
$ws.Cells.Item(27, 10).Value = 1079.1666
$ws.Cells.Item(27, 11).Value = 1732.7142
$ws.Cells.Item(27, 12).Value = 1079.1666
$ws.Cells.Item(27, 13).Value = -1625.7142
$ws.Cells.Item(27, 14).Value = -1293.1666

# LTW row 50
$ws.Cells.Item(50, 8).Value = 32333
$ws.Cells.Item(50, 9).Value = 12000
$ws.Cells.Item(50, 10).Value = 42499.5
$ws.Cells.Item(50, 11).Value = 12000
$ws.Cells.Item(50, 12).Value = 42499.5
$ws.Cells.Item(50, 13).Value = -11363
$ws.Cells.Item(50, 14).Value = -43773.5

# LTW row 55
$ws.Cells.Item(55, 8).Value = 541.9091
$ws.Cells.Item(55, 9).Value = 788
$ws.Cells.Item(55, 10).Value = 111.25
$ws.Cells.Item(55, 11).Value = 788
$ws.Cells.Item(55, 12).Value = 111.25
$ws.Cells.Item(55, 13).Value = -615
$ws.Cells.Item(55, 14).Value = -457.25

# LTW row 68
$ws.Cells.Item(68, 8).Value = 3815.2942
$ws.Cells.Item(68, 9).Value = 2766.2727
$ws.Cells.Item(68, 10).Value = 5738.5
$ws.Cells.Item(68, 11).Value = 2766.2727
$ws.Cells.Item(68, 12).Value = 5738.5
$ws.Cells.Item(68, 13).Value = -2017.2727
$ws.Cells.Item(68, 14).Value = -7236.5

# LTW row 71
$ws.Cells.Item(71, 8).Value = 3815.2942
$ws.Cells.Item(71, 9).Value = 2766.2727
$ws.Cells.Item(71, 10).Value = 5738.5
$ws.Cells.Item(71, 11).Value = 13831.3635
$ws.Cells.Item(71, 12).Value = 28692.5
$ws.Cells.Item(71, 13).Value = -10087.3635
$ws.Cells.Item(71, 14).Value = -36180.5

# LTW row 112
$ws.Cells.Item(112, 8).Value = 43949.5
$ws.Cells.Item(112, 10).Value = 43949.5
$ws.Cells.Item(112, 12).Value = 43949.5
$ws.Cells.Item(112, 14).Value = -46903.5

# LTW row 119
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 12).ClearContents()
$ws.Cells.Item(119, 14).Value = 0

# LTW row 123
$ws.Cells.Item(123, 8).Value = 94978
$ws.Cells.Item(123, 10).Value = 94978
$ws.Cells.Item(123, 12).Value = 94978
$ws.Cells.Item(123, 14).Value = -104778

$ws = $wb.Worksheets.Item("WVR")
# WVR row 105
$ws.Cells.Item(105, 8).Value = 36500
$ws.Cells.Item(105, 10).Value = 36500
$ws.Cells.Item(105, 12).Value = 36500
$ws.Cells.Item(105, 14).Value = -43488

# WVR row 113
$ws.Cells.Item(113, 8).Value = 576.25
$ws.Cells.Item(113, 9).Value = 395.375
$ws.Cells.Item(113, 10).Value = 938
$ws.Cells.Item(113, 11).Value = 1186.125
$ws.Cells.Item(113, 12).Value = 2814
$ws.Cells.Item(113, 13).Value = 983.875
$ws.Cells.Item(113, 14).Value = -7154

# WVR row 132
$ws.Cells.Item(132, 8).Value = 2633.5938
$ws.Cells.Item(132, 9).Value = 2238.8596
$ws.Cells.Item(132, 10).Value = 5847.857
$ws.Cells.Item(132, 11).Value = 6716.578799999999
$ws.Cells.Item(132, 12).Value = 17543.571
$ws.Cells.Item(132, 13).Value = -4186.578799999999
$ws.Cells.Item(132, 14).Value = -22603.571

# WVR row 136
$ws.Cells.Item(136, 8).Value = 2862.8333
$ws.Cells.Item(136, 9).Value = 1916.6511
$ws.Cells.Item(136, 10).Value = 11000
$ws.Cells.Item(136, 11).Value = 5749.9533
$ws.Cells.Item(136, 12).Value = 33000
$ws.Cells.Item(136, 13).Value = -3199.9533
$ws.Cells.Item(136, 14).Value = -38100
